$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column width adjustments (Day 3 / Day 4 narrower, Day 5 wider)
#    ColumnWidth read/write is offset ~0.83 from the stored <col width>
#    value, so subtract 0.83 from the desired stored width.
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 31.17   # F: 33 -> 32
$ws.Columns.Item(8).ColumnWidth = 31.17   # H: 33 -> 32
$ws.Columns.Item(10).ColumnWidth = 33.17  # J: 14 -> 34

# ---------------------------------------------------------------------
# 2. Activity text updates
# ---------------------------------------------------------------------
$ws.Range("J7").Value = "Master class with Ivy & Stephane"
$ws.Range("B11").Value = "Private lesson with Ivy CHUANG"
$ws.Range("H11").Value = "Private lesson with Ivy CHUANG"
$ws.Range("F20").Value = "Private lesson with Ivy CHUANG"
$ws.Range("J20").Value = "Master class with Ivy & Stephane"

# ---------------------------------------------------------------------
# 3. Drop the post-17:00 "Free Time" blocks (row 32) for Days 1-5 and
#    shrink the "Acting class" merges from 28:31 to 28:30 so they stop
#    at 17:00 instead of continuing through the removed block.
# ---------------------------------------------------------------------
$ws.Range("B28:B31").UnMerge()
$ws.Range("D28:D31").UnMerge()
$ws.Range("F28:F31").UnMerge()
$ws.Range("H28:H31").UnMerge()
$ws.Range("J28:J31").UnMerge()
$ws.Range("B32:B39").UnMerge()
$ws.Range("D32:D39").UnMerge()
$ws.Range("F32:F39").UnMerge()
$ws.Range("H32:H39").UnMerge()
$ws.Range("J32:J39").UnMerge()

# Remove the now-obsolete "Free Time" cells in row 32 entirely.
$ws.Range("B32").Clear()
$ws.Range("D32").Clear()
$ws.Range("F32").Clear()
$ws.Range("H32").Clear()
$ws.Range("J32").Clear()

# Temporarily drop the anchor cells' formatting so re-merging doesn't
# stamp the "vertical centre" style onto the blank rows 29/30 (the
# engine copies the anchor's current format across the whole merge
# range), then restore the formatting on the anchor cell afterwards.
$ws.Range("B28").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("F28").ClearFormats()
$ws.Range("H28").ClearFormats()
$ws.Range("J28").ClearFormats()

$ws.Range("B28:B30").Merge()
$ws.Range("D28:D30").Merge()
$ws.Range("F28:F30").Merge()
$ws.Range("H28:H30").Merge()
$ws.Range("J28:J30").Merge()

$ws.Range("B28").VerticalAlignment = -4108
$ws.Range("D28").VerticalAlignment = -4108
$ws.Range("F28").VerticalAlignment = -4108
$ws.Range("H28").VerticalAlignment = -4108
$ws.Range("J28").VerticalAlignment = -4108
